$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("arm")

# Remove the raw data point in J4 (Init_C for mg) - the cell becomes empty
$ws.Range("J4").ClearContents()

# Update the selection to match the new active cell
$ws.Range("J4").Select()

$wb.Application.Calculate()
